# #5: cash & deposit done
# Deposits ("存款") sheet: turn the duplicated first data row into a real
# header row, and append the standard trailing metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that the other property sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # 存款 (Deposits)

# --- Row 1: replace the old "duplicate of row 2" values with real headers
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"
$ws.Cells.Item(1,7).Value = "property_category"
$ws.Cells.Item(1,8).Value = "category"
$ws.Cells.Item(1,9).Value = "date"
$ws.Cells.Item(1,10).Value = "legislator_name"
$ws.Cells.Item(1,11).Value = "legislator_id"
$ws.Cells.Item(1,12).Value = "source_file"
$ws.Cells.Item(1,13).Value = "index"

# match the existing header-row look (bold / bordered / centered) on the
# newly added header cells
$ws.Range("B1:F1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)  # xlPasteFormats

# --- Rows 2-8: existing A-F values are untouched; append G:M metadata
$rows = @(
    @{ Row = 2;  Index = 53 },
    @{ Row = 3;  Index = 54 },
    @{ Row = 4;  Index = 55 },
    @{ Row = 5;  Index = 56 },
    @{ Row = 6;  Index = 57 },
    @{ Row = 7;  Index = 58 },
    @{ Row = 8;  Index = 59 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 7).Value  = "deposit"          # property_category
    $ws.Cells.Item($row, 8).Value  = "normal"            # category
    $ws.Cells.Item($row, 9).NumberFormat = "@"           # keep date column as text
    $ws.Cells.Item($row, 9).Value  = "2011-11-17"        # date
    $ws.Cells.Item($row, 10).Value = "蘇震清"             # legislator_name
    $ws.Cells.Item($row, 11).Value = 1718                # legislator_id
    $ws.Cells.Item($row, 12).Value = "tmp98701"          # source_file
    $ws.Cells.Item($row, 13).Value = $r.Index            # index
}

# match the existing data-row look on the newly added data cells
$ws.Range("B2:F8").Copy()
$ws.Range("G2:M8").PasteSpecial(-4122)  # xlPasteFormats
